# Edit: curvature_review_reference.docx
# 1) Remove the _GoBack bookmark from the early empty paragraph and give
#    that paragraph's mark a Cambria Math rFonts entry in its pPr/rPr.
# 2) Renumber the introduction/methods/data-analysis bookmarks (handled
#    automatically by Word -- it always renumbers w:id by document order
#    on save -- as long as the old _GoBack bookmark is gone and the new
#    bookmarks below land at the right spots).
# 3) Append a spacer paragraph, a "References" Heading1 paragraph (with a
#    "references" bookmark around the run), and two Bibliography-styled
#    reference paragraphs (Adams 2013 / Alexandre 2015), each wrapped in
#    its own "ref-*" bookmark, with a new "_GoBack" bookmark sitting in
#    the blank paragraph between them, and a trailing spacer paragraph.

$d = $word.ActiveDocument

# --- Step 1: drop the stray _GoBack bookmark on the early blank paragraph,
#     then give that (still-empty) paragraph's mark the Cambria Math font.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$earlyBlank = $d.Paragraphs.Item(4).Range
$earlyBlank = $d.Range($earlyBlank.Start, $earlyBlank.End)
$goBackParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$earlyBlank.InsertXML($goBackParaXml)

# --- Step 2: append the new block of paragraphs (spacer, References
#     heading, two bibliography entries with a blank paragraph between
#     them, trailing spacer) at the very end of the document.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$refsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
  <w:pPr>
    <w:spacing w:after="240"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>References</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Bibliography"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">Adams, D.C., and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>Otárola</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">-Castillo, E. (2013). </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>Geomorph</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">: An R package for the collection and analysis of geometric morphometric shape data. Methods in Ecology and Evolution </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:i/>
    </w:rPr>
    <w:t>4</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>, 393–399.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Bibliography"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Bibliography"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">Alexandre, H., </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>Vrignaud</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">, J., </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>Mangin</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve">, B., and Joly, S. (2015). Genetic architecture of pollination syndrome transition between hummingbird-specialist and generalist species in the genus </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:i/>
    </w:rPr>
    <w:t>Rhytidophyllum</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve"> (Gesneriaceae). </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>PeerJ</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:i/>
    </w:rPr>
    <w:t>3</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>, e1028.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="240"/>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$endRange.InsertXML($refsXml)

# --- Step 3: add the bookmarks that belong to the newly inserted
#     paragraphs. Word renumbers every bookmark's w:id by document
#     position when it saves, so we do not need to worry about the
#     order in which these (or the pre-existing) bookmarks were created.
$total = $d.Paragraphs.Count
$refPara   = $d.Paragraphs.Item($total - 4)   # "References" heading
$adamsPara = $d.Paragraphs.Item($total - 3)   # Adams, D.C. ... reference
$blankPara = $d.Paragraphs.Item($total - 2)   # blank paragraph
$alexPara  = $d.Paragraphs.Item($total - 1)   # Alexandre, H. ... reference

$refRange = $refPara.Range
$refRange = $d.Range($refRange.Start, $refRange.End)
$d.Bookmarks.Add("references", $refRange)

$adamsRange = $adamsPara.Range
$adamsRange = $d.Range($adamsRange.Start, $adamsRange.End)
$d.Bookmarks.Add("ref-adams_2013", $adamsRange)

$blankRange = $blankPara.Range
$blankRange = $d.Range($blankRange.Start, $blankRange.Start)
$d.Bookmarks.Add("_GoBack", $blankRange)

$alexRange = $alexPara.Range
$alexRange = $d.Range($alexRange.Start, $alexRange.End)
$d.Bookmarks.Add("ref-alexandre_2015", $alexRange)

Write-Output "edit complete"
